# Apply the "fM -> aM" concentration-unit rescale described in the diff:
#   - header text "Concentration (fM)" -> "Concentration (aM)"
#   - column A (rows 8:40) holds concentration values; they get multiplied by
#     1000 (fM -> aM) and the formulas (=prevCell/4) are flattened to plain
#     numeric values
#   - a new, empty column C (rows 8:40) is introduced, formatted the same way
#     (scientific notation) as column A

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header text in A1 (shared string) "Concentration (fM)" -> "Concentration (aM)"
$ws.Cells.Item(1,1).Value = "Concentration (aM)"

# 2. Column A: flatten the dilution-series formulas into static values scaled x1000
$colAValues = @{
    8  = 104000000;    9  = 26000000;     10 = 6500000;
    11 = 1625000;       12 = 406250;        13 = 101562.5;
    14 = 25390.625;     15 = 6347.65625;    16 = 1586.9140625;
    17 = 396.728515625; 18 = 99.18212890625;

    19 = 104000000;     20 = 26000000;      21 = 6500000;
    22 = 1625000;       23 = 406250;        24 = 101562.5;
    25 = 25390.625;     26 = 6347.65625;    27 = 1586.9140625;
    28 = 396.728515625; 29 = 99.18212890625;

    30 = 104000000;     31 = 26000000;      32 = 6500000;
    33 = 1625000;       34 = 406250;        35 = 101562.5;
    36 = 25390.625;     37 = 6347.65625;    38 = 1586.9140625;
    39 = 396.728515625; 40 = 99.18212890625;
}

foreach ($r in $colAValues.Keys) {
    $ws.Cells.Item($r,1).Value = $colAValues[$r]
}

# 3. New column C: empty cells, rows 8:40, same number format as column A
#    (numFmtId 11, "0.00E+00") so the existing style gets reused.
$srcFormat = $ws.Cells.Item(8,1).NumberFormat
for ($r = 8; $r -le 40; $r++) {
    $ws.Cells.Item($r,3).NumberFormat = $srcFormat
}
